$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999381899118,
    0.99728959010214235,
    0.99086583666665695,
    0.99199811655548842,
    0.97991887806468791,
    0.95088883073938324,
    0.94776064792086689,
    0.94386356805117821,
    0.94111114317478639,
    0.93929437751769618,
    0.93911823393215221,
    0.9390393395142087,
    0.92329304225947861,
    0.91912586497364979,
    0.91653445862235761,
    0.91402801814095103,
    0.91032013150307312,
    0.90921121626151136,
    0.99623637468205373,
    0.9891193370509368,
    0.98772084547317662,
    0.98645634045325969,
    0.98383148366574535,
    0.97081106234693171,
    0.96435414294761856,
    0.94723431084324838,
    0.94238898699069873,
    0.92092018230261852,
    0.90565029783787954,
    0.89908020675236955,
    0.89142665934668364,
    0.88974735376761604,
    0.88922735313091894
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
